$wb = $excel.ActiveWorkbook

# --- YDS sheet: update long simulation data strings (Week 15 added) ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = "0 3 6 2 3 2 2 -1 1 3 1 2 6 6 -4 0 21 2 0 0 3 6 4 3 11 0 4 1 7 2 5 4 4 1 0 4 4 1 -5 6 -1 25 6 4 3 46 0 4 -3 2 5 6 2 9 4 4 4 -1 2 1 0 2 13 4 5 4 9 4 3 5 7 8 1 0 7 23 4 2 1 2 3 3 2 4 -1 10 5 0 2 14 4 4 11 4 2 6 0 5 3 3 6 37 10 -1 2 -2 2 1 4 2 5 2 3 2 4 3 11 2 1 4 2 2 25 -3 1 4 14 5 7 3 5 3 12 7 2 2 1 19 0 -1 13 8 9 -1 1 3 1 24 3 -1 6 6 -3 -3 2 3 1 2 1 -1 7 2 4 5 5 2 13 1 2 2 7 2 3 7 7 0 3 3 2 -2 2 7 0 7 2 0 0 1 8 3 0 7 0 3 1 14 4 2 -3 1 6 3 11 1 98 1 -3 3 -2 1 -1 4 5 4 12 1 6 6 2 3 8 0 2 1 2 6 5 1 2 -1 1 2 1 5 2 6 3 2 5 4 1 2 11 3 1 34 1 6 5 2 5 3 1 5 3 3 3 1 1 10 5 3 14 4 2 16 3 -3 7 14 1 2 1 3 1 9 4 1 0 1 1 17 2 2 6 2 1 2 8 4 2 3 4 6 4 14 12 9 2 5 3 1 9 1 2 -1 1 4 4 -4 3 6 2 5 16 7 9 2 -3 0 4 9 8 3 5 6 3 8 9 7 7 -2 6 0 4 3 22 3 8 4 3 2 3 3 6 6 5 17 8 2 3 2 13 1 2 4 9 0 5 8 2 6 3 7 2 3 7 3 3 10 5 6 4 1 12 1 2 3 1 0 2 3 4 5 3 3 1 3 4 2 0 1 20 1 5 9 2 5 2 3 -1 2 2 4 -1 1 6 3 13 -3 3 6 6 11 7 8 4 1 1 0 3 4 27 6 3 2 0 4 7 8 4 4 8 1 1 5 5 3 5 3 3 5 2 5 1 7 0 6 6 1 8 2 7 3 2 10 3 8 10 4 2 8 4 6 0 3 1 3 1 5 -1 4 3 0 2 3 1 4 9 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 1 17 4 8 4 4 4 6 4 -1 13 13 3 6 1 4 5 5 4 6 2 4 0 9 4 3 6 0 2 3 2 2 20 5 6 2 7 0 1 9 2 2 0 1 -1 3 15 2 4 5 12 3 9 4 11 6 8 12 2 -1 1 12 19 12 9 5 3 1 2 -3 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 5 6 0 3 5 8 2 1 3 4 8 6 7 10 7 6 3 2 5 10 3 3 4 6 2 7 -1 3 4 9 2 2 1 3 -3 2 -1 4 7 3 7 -3 4 2 2 1 11 6 12 4 1 4 5 8 8 1 2 3 1 11 8 3 28 3 4 2 2 6 6 5 4 3 1 3 4 9 -4 1 9 0 3 47 4 2 12 2 3 1 5 -1 8 2 1 1 2 9 1 -2 2 0 7 21 1 4 3 1"
$ydsWs.Range("C2").Value = "4 9 2 1 2 13 7 4 -4 -4 6 4 1 8 3 4 -2 -3 5 3 -5 4 4 2 4 0 -3 2 4 4 -2 5 9 1 4 3 -1 3 7 1 5 -4 5 4 8 5 1 -1 4 5 14 1 2 9 2 -2 -1 2 3 9 8 0 2 2 10 2 4 5 0 3 -1 2 4 3 -1 3 -3 2 2 1 0 0 3 4 2 3 2 2 4 3 3 1 4 7 3 1 1 0 2 5 -2 25 0 6 1 0 7 4 1 7 2 -3 8 3 20 4 3 1 3 0 5 4 -1 2 1 10 0 4 1 2 2 9 10 2 4 0 9 2 1 6 2 -2 2 2 1 9 1 6 11 3 17 6 4 2 2 1 2 6 5 4 -1 15 1 0 2 17 2 1 23 5 12 1 11 -2 4 10 -1 5 17 -2 1 5 -1 3 3 7 -1 4 1 0 6 5 5 6 4 7 6 1 7 6 -1 4 1 10 -3 8 1 3 4 2 1 1 -2 -1 0 11 -1 2 0 2 -5 14 2 0 6 -3 2 1 0 8 -2 2 3 4 10 10 17 5 0 0 5 4 7 4 4 8 1 5 13 7 2 5 3 8 2 5 1 3 4 9 4 4 1 12 3 14 0 3 4 0 4 0 4 4 10 7 1 12 7 -1 5 1 16 -4 3 3 0 3 2 7 -2 4 7 14 4 9 4 -4 3 3 0 0 1 16 0 11 3 18 1 5 0 5 0 4 0 4 2 1 4 0 1 2 2 -1 -6 9 3 11 4 9 62 1 2 4 -4 4 -1 2 10 2 2 5 2 11 4 1 10 3 3 3 3 -1 13 3 -1 8 4 0 2 3 10 5 1 1 11 6 5 1 5 2 5 8 1 2 17 6 1 2 2 2 4 3 1 2 12 12 1 7 4 0 6 1 5 3 1 0 9 11 2 5 5 11 6 -1 2 10 7 26 10 3 3 7 1 4 2 4 0 2 7 1 0 3 -1 13 3 0 2 4 1 4 11 7 -1 0 2 3 10 0 0 8 2 0 8 1 -3 1 8 6 0 3 6 4 0 4 4 2 5 1 4 -2 0 2 11 2 3 1 4 1 5 15 0 6 -4 0 0 0 -1 1 2 4 0 3 1 13 7 4 5 -2 11 1 1 -1 5 6 2 5 -2 9 11 23 14 7 3 -1 2 0 29 1 1 6 4 0 6 13 5 12 6 0 6 0 8 8 1 2 -1 4 -1 -1 1 1 8 2 6 16 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 1 3 1 5 -1 -8 10 5 2 3 15 4 4 -5 -1 -3 12 9 1 1 1 2 3 4 7 3 2 4 2 -1 5 0 3 2 8 6 0 1 3 16 1 2 12 0 9 6 0 4 -2 0 5 12 5 1 16 3 2 5 15 5 10 15 1 3 4 7 39 17 6 5 1 5 11 0 4 2 0 4 0 3 1 0 2 3 0 14 0 14 -3 10 22 1 29 18 4 11 5 18 23 2 0 7 8 3"
$ydsWs.Range("B3").Value = "29 8 6 7 2 10 6 14 9 12 9 23 9 15 5 37 2 3 2 6 7 2 50 6 23 2 7 4 6 12 8 36 12 11 2 0 5 8 -1 10 5 9 7 3 9 10 10 47 6 6 1 33 9 1 27 6 17 11 6 7 10 8 11 7 30 9 26 3 7 24 17 3 6 -2 15 5 1 5 10 6 9 5 10 12 28 44 19 6 7 48 9 7 13 29 11 17 2 35 9 2 0 9 4 9 10 1 8 6 14 1 11 7 23 9 16 12 4 6 9 10 13 2 6 4 6 0 7 8 25 3 12 31 6 18 3 6 22 12 8 8 6 11 16 26 5 5 2 16 11 7 33 6 35 3 8 10 2 4 13 15 6 12 4 23 22 1 9 19 3 10 3 7 14 3 9 2 8 9 15 12 8 1 9 8 4 16 3 20 25 8 18 5 4 7 11 20 7 6 2 2 10 22 0 4 15 2 16 14 6 20 14 6 23 19 14 3 5 3 9 8 7 2 4 15 16 44 3 6 6 19 22 6 8 3 25 13 7 6 7 6 7 16 4 9 4 7 9 18 12 4 7 1 10 9 9 9 10 10 13 7 12 12 6 6 2 6 8 48 -4 44 9 14 9 8 7 5 7 8 6 10 11 7 9 20 48 16 6 4 3 9 7 29 2 11 10 6 6 14 -8 14 8 2 7 10 7 10 20 14 32 5 26 16 24 4 17 8 5 11 9 13 4 18 13 46 13 8 33 33 47 -2 27 7 16 10 12 12 13 7 20 22 7 7 15 11 16 12 25 5 30 7 4 5 10 12 11 29 5 14 23 20 25 12 8 46 7 8 12 14 18 47 4 30 3 6 2 15 17 12 15 36 4 5 27 22 19 24 13 13 13 7 23 17 8 20 19 35 16 5 11 16 3 8 8 14 8 10 10 8 2 19 6 15 14 13 27 14 15 5 -3 14 52 7 19 1 6 39 8 16 -1 6 11 29 9 1 5 16 15 5 8 3 31 -3 9 4 17 8 5 15 1 3 12 25 15 6 10 19 28 16 8 5 27 5 19 17 2 10 47 20 2 2 11 10 1 5 7 9 7 9 7 3 11 8 4 6 20 24 20 25 20 16 16 9 5 8 1 9 5 26 17 21 3 8 1 2 9 10 26 5 12 6 8 5 1 15 17 4 3 7 7 0 17 5 26 3 20 20 32 14 8 20 15 11 6 7 17 11 -4 23 13 31 17 8 -2 11 8 5 5 21 -1 9 9 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 9 26 13 10 10 18 6 15 62 6 10 10 11 12 3 4 11 5 12 6 14 12 16 4 15 7 34 19 9 22 9 22 16 10 2 8 8 11 3 9 23 -4 10 12 14 18 10 3 6 -4 5 4 4 9 15 0 8 5 4 8 10 2 27 10 -4 6 4 21 4 11 9 10 46 2 9 14 0 9 8 9 7 12 0 34 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 2 3 29 9 8 3 6 7 11 4 5 14 15 8 6 2 15 4 0 4 3 40 10 8 15 3 13 7 10 8 6 6 6 7 4 35 7 12 16 3 13 9 17 14 14 18 4 5 8 3 21 5 3 20 9 0 0 4 -1 10 8 7 7 15 9 26 5 4 11 10 32 9 2 9 3 8 6 13 15 5 11 0 6 3 15 7 4 3 25 36 3 0 15 6 17 12 27 1 9 1 7 11 5 7 18 1 5 11 14 6 6 6 36 5 3 6 15 5 7 5 9 20 8 16 23 18 5 13 16 13 7 5 20 5 6 24 4 18 0 8 5 11 10 6 -7 14 58"
$ydsWs.Range("C3").Value = "14 -2 5 5 1 4 -2 18 10 12 4 17 8 46 6 5 5 4 38 6 10 19 -1 8 9 7 13 12 3 3 7 4 15 5 10 11 30 39 16 18 8 7 5 12 12 2 22 14 8 13 7 14 9 6 2 15 18 7 -2 4 15 28 5 9 5 4 4 5 2 22 7 6 3 2 10 53 3 3 5 9 1 17 19 7 17 8 29 10 9 1 7 3 72 8 3 10 5 24 13 8 7 9 -6 4 12 13 25 7 8 12 3 9 8 9 8 6 8 5 2 10 -1 17 4 20 8 12 12 11 15 16 7 9 5 5 18 4 -1 6 28 14 20 5 7 29 9 6 11 21 6 9 15 6 11 2 44 1 11 3 10 6 3 7 8 18 0 7 11 5 12 14 11 11 12 6 15 4 11 9 1 12 1 6 10 19 6 8 11 6 19 14 40 5 2 11 8 7 -6 21 12 9 20 12 3 5 5 10 -5 5 29 9 1 3 12 5 0 2 5 -3 7 6 38 24 7 5 1 5 17 2 2 4 17 5 6 8 5 37 7 4 -1 6 3 8 19 7 24 6 7 5 9 8 4 35 -1 12 6 8 16 2 7 18 18 4 5 17 5 -3 -1 23 9 19 34 23 6 3 1 75 19 3 14 7 44 20 12 14 2 9 5 13 9 2 8 29 37 7 7 8 9 14 8 2 2 4 6 20 8 11 3 15 3 2 8 3 15 15 2 6 17 4 5 15 15 14 8 9 5 5 14 25 9 2 4 5 2 10 4 23 27 5 3 25 5 5 8 13 14 12 2 7 14 18 5 15 12 -1 13 28 4 14 7 4 -4 11 33 -2 15 2 7 2 2 7 4 19 11 2 6 35 22 15 3 10 6 5 20 2 9 7 7 6 12 7 26 4 13 13 7 4 8 2 11 10 16 6 7 5 15 6 17 11 19 14 13 0 3 11 18 -1 5 36 24 18 10 5 5 8 12 7 12 29 12 19 6 1 8 23 4 3 5 1 16 11 16 10 7 -1 13 4 4 56 10 16 -3 9 7 13 16 0 8 8 6 7 6 14 0 2 23 15 5 50 12 4 5 9 5 4 23 5 2 10 24 10 16 8 2 4 13 7 5 11 7 2 5 29 9 11 3 6 14 13 4 12 11 -1 5 8 0 5 21 13 17 33 12 9 4 23 8 18 16 1 4 28 5 16 11 8 12 7 8 7 22 2 5 22 6 8 3 4 13 -5 4 11 7 10 5 16 6 14 6 4 3 7 10 9 6 21 2 3 6 11 9 31 10 0 9 5 9 4 5 23 24 4 8 24 7 6 1 14 8 3 8 14 4 23 12 8 7 22 5 7 8 1 2 4 10 2 6 7 6 12 22 11 5 6 6 6 9 2 11 21 22 2 75 0 20 5 40 10 5 17 6 4 10 3 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 21 4 34 9 12 10 5 24 5 3 12 5 2 21 6 10 8 4 23 10 2 7 1 0 11 5 21 23 5 6 4 16 2 1 6 -3 24 25 6 6 5 5 5 3 6 3 17 4 12 4 10 22 20 6 22 12 -6 11 8 5 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 8 14 13 4 8 11 1 20 22 7 9 8 14 1 15 6 18 21 9 6 16 6 5 6 6 9 12 9 4 3 5 5 0 11 21 4 28 8 11 8 -1 7 2 4 5 5 8 4 9 8 62 10 14 3 3 15 9 5 12 22 9 23 -3 15 5 15 11 11 7 8 7 10 0 6 6 12 19 29 7 10 9 8 9 4 -3 14 2 19 6 9 4 16 6 20 6 17 8 8 8 2 7 30 1 10 11 5 3 25 24 21 9 5 6 10 12 6 1 13 9 -2 10 2 4 10 15 9 12 12 3 4 4 4 2 15 1 24 3"

# --- OFF sheet: updated aggregate totals ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 409
$offWs.Range("D2").Value = 24
$offWs.Range("F2").Value = 73
$offWs.Range("G2").Value = 90
$offWs.Range("J2").Value = 50
$offWs.Range("L2").Value = 557
$offWs.Range("M2").Value = 358
$offWs.Range("O2").Value = 30
$offWs.Range("Q2").Value = 881
$offWs.Range("C3").Value = 371
$offWs.Range("D3").Value = 17
$offWs.Range("E3").Value = 60
$offWs.Range("F3").Value = 291
$offWs.Range("G3").Value = 95
$offWs.Range("H3").Value = 36
$offWs.Range("I3").Value = 130
$offWs.Range("J3").Value = 128
$offWs.Range("N3").Value = 14

# --- DEF sheet: updated aggregate totals ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 346
$defWs.Range("D2").Value = 15
$defWs.Range("E2").Value = 21
$defWs.Range("F2").Value = 86
$defWs.Range("G2").Value = 90
$defWs.Range("J2").Value = 50
$defWs.Range("L2").Value = 645
$defWs.Range("M2").Value = 444
$defWs.Range("O2").Value = 47
$defWs.Range("P2").Value = 26
$defWs.Range("Q2").Value = 1006
$defWs.Range("B3").Value = 17
$defWs.Range("C3").Value = 453
$defWs.Range("D3").Value = 16
$defWs.Range("E3").Value = 72
$defWs.Range("F3").Value = 262
$defWs.Range("G3").Value = 96
$defWs.Range("H3").Value = 64
$defWs.Range("I3").Value = 119
$defWs.Range("J3").Value = 132
$defWs.Range("N3").Value = 40

# --- ST sheet: updated aggregate totals + simulation data strings ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 202
$stWs.Range("D2").Value = 107
$stWs.Range("F2").Value = 455
$stWs.Range("G2").Value = 437
$stWs.Range("J2").Value = 183
$stWs.Range("K2").Value = 173
$stWs.Range("B3").Value = 170
$stWs.Range("B4").Value = "73 74 74 65 59 66 46 68 67 72 72 65 60 66 61 69 67 65 51 50 44 56 76 64 67 74 73 64 64 69 67 69"
$stWs.Range("B5").Value = "32 47 40 44 23 31 0 28 98 25 27 27 26 31 30 41 21 25 13 19 14 45 27 0 25 33 36 28 19 72 30 29"
$stWs.Range("B6").Value = "27 33 24 25 20 21 15 28 34 28 13 25 26 21 25 24 6 23 22 22 13 30 30 21 28 23 18 8 29 22 18 43 9 26 28 21 20 41 31 16 24 19 24 24 13 10 17 17 17 19 34 18 13 15 24 30 12 15 33 7"
$stWs.Range("D3").Value = "43 50 38 59 46 40 47 32 36 41 42 42 58 56 56 40 44 46 30 47 25 44 43 49 36 44 49 33 33 49 39 40 59 32 56 46 44 49 62 53 61 43 47 39 52 40 41 37 39 49 44 44 53 30 43 48 48 44 46 34 40 33 38 39 65 42 49 41 48 41 43 57 38 53 39 15 55 38 34 43 41 41 44 42 49 42 46 54 32 40 43 46 35 41 55 40 14 54 55 60 41 32 55 20 39 35 47"
$stWs.Range("D4").Value = "0 18 16 18 7 0 -1 0 0 0 0 0 0 0 0 0 6 0 0 0 0 0 0 0 0 0 11 0 0 7 8 15 0 0 2 0 12 0 5 3 10 5 0 0 4 0 0 0 0 74 4 0 0 0 54 7 0 0 0 0 0 0 0 0 0 3 12 0 6 0 5 0 9 0 0 0 0 0 11 6 0 0 0 8 0 13 0 0 0 12 10 0 0 0 0 0 0 9 13 7 0 0 4 0 0 0 5"
$stWs.Range("D5").Value = "0 0 10 0 0 0 0 0"

# --- TURNS sheet: updated aggregate totals ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C2").Value = 19
$turnsWs.Range("D2").Value = 8
$turnsWs.Range("E3").Value = 19

# --- PEN sheet: updated aggregate totals ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D2").Value = 39
